# Updates the crypto price/volume table on Sheet1 to the latest scraped
# values (GitHub Actions refresh). Column D ("Price") cells whose new
# value is purely numeric-looking are written with a leading apostrophe
# so Excel keeps them as text (matching the original inline-string/
# general-format cells) instead of silently parsing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.916.45"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.805.61"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'310.07"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4409"
$ws.Range("E7").Value = "  +4.35%  "
$ws.Range("D8").Value = "'0.3708"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "'0.07447"
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").Value = "'0.8609"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.796.04"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "'6.646"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'93.15"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "'0.07070"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'5.279"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'0.000008697"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'0.9998"
$ws.Range("D20").Value = "'14.83"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "26.946.64"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'5.172"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "'10.83"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "2.019.68"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("D25").Value = "'1.984"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "'151.28"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "'18.37"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'117.43"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'0.08775"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'0.7447"
$ws.Range("D33").Value = "'1.165"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'4.482"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "'2.888"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("D38").Value = "'0.01973"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("D39").Value = "'0.05212"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "'0.5252"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").Value = "'7.077"
$ws.Range("D42").Value = "'2.821"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Value = "'0.1685"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.123"
$ws.Range("E44").Value = "  +9.81%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.499"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'0.4978"
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("D47").Value = "'10.38"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "'104.23"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("D49").Value = "'0.9995"
$ws.Range("D50").Value = "'1.668"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "'0.06345"
$ws.Range("E51").Value = "  +0.37%  "
